$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.171.51'
$ws.Range("E2").Value = '  -4.04%  '

$ws.Range("D3").Value = '1.658.14'
$ws.Range("E3").Value = '  -2.66%  '

$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").Value = '''217.79'
$ws.Range("E5").Value = '  -2.51%  '

$ws.Range("D6").Value = '''0.5141'
$ws.Range("E6").Value = '  -3.09%  '

$ws.Range("E7").Value = '  +0.25%  '

$ws.Range("D8").Value = '''0.2581'
$ws.Range("E8").Value = '  -2.88%  '

$ws.Range("D9").Value = '''0.06432'
$ws.Range("E9").Value = '  -2.20%  '

$ws.Range("D10").Value = '''19.95'
$ws.Range("E10").Value = '  -3.67%  '

$ws.Range("D11").Value = '''0.07805'
$ws.Range("E11").Value = '  +2.49%  '

$ws.Range("D12").Value = '1.668.37'
$ws.Range("E12").Value = '  -2.29%  '

$ws.Range("D13").Value = '''4.294'
$ws.Range("E13").Value = '  -4.34%  '

$ws.Range("D14").Value = '1.885.38'
$ws.Range("E14").Value = '  -2.79%  '

$ws.Range("D15").Value = '''0.5547'
$ws.Range("E15").Value = '  -3.96%  '

$ws.Range("D16").Value = '0.0₅8059'
$ws.Range("E16").Value = '  -0.98%  '

$ws.Range("D17").Value = '''64.23'
$ws.Range("E17").Value = '  -4.80%  '

$ws.Range("D18").Value = '26.217.12'
$ws.Range("E18").Value = '  -3.93%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '''1.006'
$ws.Range("E19").Value = '  +0.36%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''211.62'
$ws.Range("E20").Value = '  -1.57%  '

$ws.Range("D21").Value = '''4.431'
$ws.Range("E21").Value = '  -3.79%  '

$ws.Range("D22").Value = '''10.04'
$ws.Range("E22").Value = '  -2.96%  '

$ws.Range("D23").Value = '''5.976'
$ws.Range("E23").Value = '  +0.29%  '

$ws.Range("E24").Value = '  +0.25%  '

$ws.Range("D25").Value = '''143.58'
$ws.Range("E25").Value = '  -0.20%  '

$ws.Range("D26").Value = '''1.756'
$ws.Range("E26").Value = '  +3.30%  '

$ws.Range("D27").Value = '''0.1165'
$ws.Range("E27").Value = '  -2.76%  '

$ws.Range("D28").Value = '''6.972'
$ws.Range("E28").Value = '  -3.12%  '

$ws.Range("D29").Value = '''15.79'

$ws.Range("D30").Value = '''0.05219'
$ws.Range("E30").Value = '  -2.62%  '

$ws.Range("E31").Value = '  -2.46%  '

$ws.Range("D32").Value = '''3.367'
$ws.Range("E32").Value = '  -2.70%  '

$ws.Range("D33").Value = '''3.218'
$ws.Range("E33").Value = '  -5.31%  '

$ws.Range("D34").Value = '''1.569'
$ws.Range("E34").Value = '  -4.26%  '

$ws.Range("D35").Value = '''2.761'
$ws.Range("E35").Value = '  -3.56%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''0.9307'
$ws.Range("E36").Value = '  -1.51%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '''2.371'
$ws.Range("E37").Value = '  -1.72%  '

$ws.Range("D38").Value = '1.172.40'
$ws.Range("E38").Value = '  +12.75%  '

$ws.Range("D39").Value = '''0.5698'
$ws.Range("E39").Value = '  -1.58%  '

$ws.Range("D40").Value = '''0.01594'
$ws.Range("E40").Value = '  -1.82%  '

$ws.Range("E41").Value = '  +0.25%  '

$ws.Range("D42").Value = '''0.8436'
$ws.Range("E42").Value = '  +0.47%  '

$ws.Range("D43").Value = '''5.670'
$ws.Range("E43").Value = '  -1.88%  '

$ws.Range("D44").Value = '''100.52'
$ws.Range("E44").Value = '  -0.41%  '

$ws.Range("D45").Value = '1.795.82'
$ws.Range("E45").Value = '  -2.80%  '

$ws.Range("D46").Value = '0.0₈114'
$ws.Range("E46").Value = '  +1.25%  '

$ws.Range("D47").Value = '''0.4537'
$ws.Range("E47").Value = '  +0.52%  '

$ws.Range("D48").Value = '''55.91'
$ws.Range("E48").Value = '  -3.12%  '

$ws.Range("D49").Value = '''1.003'
$ws.Range("E49").Value = '  -0.22%  '

$ws.Range("D50").Value = '''7.847'
$ws.Range("E50").Value = '  -2.68%  '

$ws.Range("D51").Value = '''0.05056'
